$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the FechaSiniestro date string (shared string used by G2).
# Setting .Value alone resets the cell's quote-prefix style (s="1"), so
# re-apply formats-only from a sibling cell that carries the same style
# (G3) after writing the new value.
$ws.Range("G2").Value = "19/03/2021"
$ws.Range("G3").Copy()
$ws.Range("G2").PasteSpecial(-4122)  # xlPasteFormats

# Update NroPoliza value in E2, likewise re-applying its quote-prefix
# style (s="1") from F2 afterwards.
$ws.Range("E2").Value = 11111003014
$ws.Range("F2").Copy()
$ws.Range("E2").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# Change active selection to G2
$ws.Range("G2").Select()
